# Populate "Sheet1" with the MOC (margin-of-confidence / curve) data set.
# Column A = x values, Column B = y values, 21 rows (A1:B21).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0, 30),
    @(53.63194324576201, 59.613313518171438),
    @(78.486691304574578, 72.530028314322806),
    @(94.918730416448554, 80.294763843209921),
    @(111.79351904740416, 87.754027900942958),
    @(129.7956344149467, 95.174652557836396),
    @(149.37001427553867, 102.67154217157351),
    @(170.8932301192836, 110.29809349679101),
    @(194.73311864238482, 118.07452690361883),
    @(221.27723852371085, 125.99813816483909),
    @(250.9507883081364, 134.04661829440482),
    @(284.23093220033189, 142.17801678043591),
    @(321.66042839806977, 150.32863238924617),
    @(363.86205656024242, 158.40926842420325),
    @(411.55483002827003, 166.29990987215484),
    @(465.57281656215304, 173.84266621131829),
    @(526.88739348644356, 180.83266592418224),
    @(596.63386204558617, 187.00643948272096),
    @(676.14351949012485, 192.02716081773431),
    @(766.98253451320784, 195.46591381996029),
    @(870.99930250701914, 196.77789127236366)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $ws.Cells.Item($i + 1, 1).Value = $row[0]
    $ws.Cells.Item($i + 1, 2).Value = $row[1]
}

# Match the saved file's selection (A1:B21 highlighted, default active cell).
$null = $ws.Range("A1:B21").Select()
